$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.920.95'
$ws.Range("D3").Value = '''1.810.23'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''310.31'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '''0.4619'
$ws.Range("E7").Value = '  +3.64%  '
$ws.Range("D8").Value = '''0.3712'
$ws.Range("E8").Value = '  -1.62%  '
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("D10").Value = '''0.8759'
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("E11").Value = '  -1.99%  '
$ws.Range("D12").Value = '''1.885.67'
$ws.Range("E12").Value = '  +3.12%  '
$ws.Range("D13").Value = '''5.363'
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = '''92.26'
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").Value = '''6.508'
$ws.Range("E15").Value = '  -3.25%  '
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '''0.000008706'
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '''14.75'
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").Value = '''26.907.67'
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").Value = '''5.326'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("E23").Value = '  -2.93%  '
$ws.Range("D24").Value = '''2.016.97'
$ws.Range("E24").Value = '  -2.10%  '
$ws.Range("D25").Value = '''1.895'
$ws.Range("E25").Value = '  -2.65%  '
$ws.Range("D26").Value = '''151.34'
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  -1.19%  '
$ws.Range("D28").Value = '''2.155'
$ws.Range("E28").Value = '  -5.92%  '
$ws.Range("D29").Value = '''5.346'
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("D30").Value = '''115.99'
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("D31").Value = '''0.08892'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = '''0.7573'
$ws.Range("E32").Value = '  -4.95%  '
$ws.Range("D33").Value = '''1.158'
$ws.Range("E33").Value = '  -3.59%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''2.917'
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '''4.454'
$ws.Range("E35").Value = '  -2.43%  '
$ws.Range("D36").Value = '''1.000'
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").Value = '''1.103'
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("D38").Value = '''0.01969'
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("D39").Value = '''0.05248'
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("D40").Value = '''2.427'
$ws.Range("E40").Value = '  +3.09%  '
$ws.Range("D41").Value = '''2.932'
$ws.Range("E41").Value = '  +2.04%  '
$ws.Range("D42").Value = '''0.5328'
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").Value = '''7.204'
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("E44").Value = '  -2.33%  '
$ws.Range("D45").Value = '''8.509'
$ws.Range("E45").Value = '  -2.28%  '
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("D47").Value = '''10.32'
$ws.Range("E47").Value = '  -2.85%  '
$ws.Range("D48").Value = '''1.000'
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("D50").Value = '''103.84'
$ws.Range("E50").Value = '  -1.67%  '
$ws.Range("D51").Value = '''0.06297'
$ws.Range("E51").Value = '  -1.51%  '
